# Apply the price/volume refresh from the "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.709.09'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.597.68'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.26'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0620'
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.53'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '1.822.25'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').Value = '1.593.60'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('E15').Value = '  +0.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.39'
$ws.Range('E16').Value = '  +1.47%  '
$ws.Range('D17').Value = '0.0₃0766'
$ws.Range('E17').Value = '  +5.16%  '
$ws.Range('D18').Value = '26.671.52'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '209.74'
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.00'
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.19'
$ws.Range('E21').Value = '  +5.63%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.95'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.13'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.15'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.32'
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0518'
$ws.Range('E30').Value = '  +3.00%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('D34').Value = '1.288.16'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  -5.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.46'
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  +17.43%  '
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.785'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.19'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.27'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').Value = '1.734.95'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.32'
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.57'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.101'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.33'
$ws.Range('E51').Value = '  -1.32%  '
